$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3117.879
$ws.Range("I113").Value = 2430
$ws.Range("J113").Value = 3565
$ws.Range("K113").Value = 2430
$ws.Range("L113").Value = 3565
$ws.Range("M113").Value = 824
$ws.Range("N113").Value = -10073
$ws.Range("H135").Value = 488.58334
$ws.Range("I135").Value = 510.21875
$ws.Range("J135").Value = 445.3125
$ws.Range("K135").Value = 4591.96875
$ws.Range("L135").Value = 4007.8125
$ws.Range("M135").Value = -2056.96875
$ws.Range("N135").Value = -9077.8125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2082.625
$ws.Range("I2").Value = 2944.6667
$ws.Range("J2").Value = 974.2857
$ws.Range("K2").Value = 2944.6667
$ws.Range("L2").Value = 974.2857
$ws.Range("M2").Value = -2831.6667
$ws.Range("N2").Value = -1200.2857
$ws.Range("H74").Value = 9823320
$ws.Range("I74").Value = 6867964
$ws.Range("J74").Value = 17600570
$ws.Range("K74").Value = 6867964
$ws.Range("L74").Value = 17600570
$ws.Range("M74").Value = -6867090
$ws.Range("N74").Value = -17602318
$ws.Range("H77").Value = 9823320
$ws.Range("I77").Value = 6867964
$ws.Range("J77").Value = 17600570
$ws.Range("K77").Value = 34339820
$ws.Range("L77").Value = 88002850
$ws.Range("M77").Value = -34335452
$ws.Range("N77").Value = -88011586
$ws.Range("H97").Value = 1465.6
$ws.Range("I97").Value = 1459.5238
$ws.Range("J97").Value = 1497.5
$ws.Range("K97").Value = 1459.5238
$ws.Range("L97").Value = 1497.5
$ws.Range("M97").Value = -963.5237999999999
$ws.Range("N97").Value = -2489.5
$ws.Range("H116").Value = 2082.625
$ws.Range("I116").Value = 2944.6667
$ws.Range("J116").Value = 974.2857
$ws.Range("K116").Value = 2944.6667
$ws.Range("L116").Value = 974.2857
$ws.Range("M116").Value = -650.6667000000002
$ws.Range("N116").Value = -5562.2857
$ws.Range("H122").Value = 2110.7693
$ws.Range("I122").Value = 2275.2122
$ws.Range("J122").Value = 1206.3334
$ws.Range("K122").Value = 6825.6366
$ws.Range("L122").Value = 3619.0002
$ws.Range("M122").Value = -4375.6366
$ws.Range("N122").Value = -8519.0002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2082.625
$ws.Range("I3").Value = 2944.6667
$ws.Range("J3").Value = 974.2857
$ws.Range("K3").Value = 2944.6667
$ws.Range("L3").Value = 974.2857
$ws.Range("M3").Value = -2830.6667
$ws.Range("N3").Value = -1202.2857
$ws.Range("H64").Value = 494.55554
$ws.Range("I64").Value = 290.5
$ws.Range("J64").Value = 552.8570999999999
$ws.Range("K64").Value = 290.5
$ws.Range("L64").Value = 552.8570999999999
$ws.Range("M64").Value = -65.5
$ws.Range("N64").Value = -1002.8571
$ws.Range("H67").Value = 494.55554
$ws.Range("I67").Value = 290.5
$ws.Range("J67").Value = 552.8570999999999
$ws.Range("K67").Value = 290.5
$ws.Range("L67").Value = 552.8570999999999
$ws.Range("M67").Value = 489.5
$ws.Range("N67").Value = -2112.8571
$ws.Range("H80").Value = 667.63635
$ws.Range("I80").Value = 298
$ws.Range("J80").Value = 704.6
$ws.Range("K80").Value = 298
$ws.Range("L80").Value = 704.6
$ws.Range("M80").Value = 700
$ws.Range("N80").Value = -2700.6
$ws.Range("H83").Value = 667.63635
$ws.Range("I83").Value = 298
$ws.Range("J83").Value = 704.6
$ws.Range("K83").Value = 1490
$ws.Range("L83").Value = 3523
$ws.Range("M83").Value = 3502
$ws.Range("N83").Value = -13507
$ws.Range("H107").Value = 1488.8108
$ws.Range("I107").Value = 1528.3793
$ws.Range("K107").Value = 1528.3793
$ws.Range("M107").Value = 391.6206999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 84009.164
$ws.Range("I16").Value = 100629.7
$ws.Range("J16").Value = 906.5
$ws.Range("K16").Value = 100629.7
$ws.Range("L16").Value = 906.5
$ws.Range("M16").Value = -100342.7
$ws.Range("N16").Value = -1480.5
$ws.Range("H58").Value = 5122.2104
$ws.Range("I58").Value = 5697.76
$ws.Range("J58").Value = 4015.3845
$ws.Range("K58").Value = 5697.76
$ws.Range("L58").Value = 4015.3845
$ws.Range("M58").Value = -5494.76
$ws.Range("N58").Value = -4421.3845
$ws.Range("H94").Value = 6985.7144
$ws.Range("I94").Value = 1200
$ws.Range("J94").Value = 9300
$ws.Range("K94").Value = 1200
$ws.Range("L94").Value = 9300
$ws.Range("M94").Value = -749
$ws.Range("N94").Value = -10202
$ws.Range("H99").Value = 47859.41
$ws.Range("I99").Value = 73558.78999999999
$ws.Range("J99").Value = 2885.5
$ws.Range("K99").Value = 73558.78999999999
$ws.Range("L99").Value = 2885.5
$ws.Range("M99").Value = -72060.78999999999
$ws.Range("N99").Value = -5881.5
$ws.Range("H107").Value = 251.7561
$ws.Range("I107").Value = 190.28
$ws.Range("J107").Value = 347.8125
$ws.Range("K107").Value = 190.28
$ws.Range("L107").Value = 347.8125
$ws.Range("M107").Value = 1729.72
$ws.Range("N107").Value = -4187.8125
$ws.Range("H113").Value = 84009.164
$ws.Range("I113").Value = 100629.7
$ws.Range("J113").Value = 906.5
$ws.Range("K113").Value = 100629.7
$ws.Range("L113").Value = 906.5
$ws.Range("M113").Value = -98459.7
$ws.Range("N113").Value = -5246.5
$ws.Range("H126").Value = 47859.41
$ws.Range("I126").Value = 73558.78999999999
$ws.Range("J126").Value = 2885.5
$ws.Range("K126").Value = 220676.37
$ws.Range("L126").Value = 8656.5
$ws.Range("M126").Value = -218206.37
$ws.Range("N126").Value = -13596.5
$ws.Range("H132").Value = 25003462
$ws.Range("I132").Value = 71429490
$ws.Range("J132").Value = 4830.615
$ws.Range("K132").Value = 214288470
$ws.Range("L132").Value = 14491.845
$ws.Range("M132").Value = -214285940
$ws.Range("N132").Value = -19551.845
$ws.Range("H136").Value = 5122.2104
$ws.Range("I136").Value = 5697.76
$ws.Range("J136").Value = 4015.3845
$ws.Range("K136").Value = 17093.28
$ws.Range("L136").Value = 12046.1535
$ws.Range("M136").Value = -14543.28
$ws.Range("N136").Value = -17146.1535

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 15013.5
$ws.Range("I49").Value = 10027
$ws.Range("J49").Value = 20000
$ws.Range("K49").Value = 10027
$ws.Range("L49").Value = 20000
$ws.Range("M49").Value = -9843
$ws.Range("N49").Value = -20368
$ws.Range("H113").Value = 909.2353000000001
$ws.Range("I113").Value = 715.5454999999999
$ws.Range("J113").Value = 1264.3334
$ws.Range("K113").Value = 715.5454999999999
$ws.Range("L113").Value = 1264.3334
$ws.Range("M113").Value = 1454.4545
$ws.Range("N113").Value = -5604.3334

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 16532.572
$ws.Range("J42").Value = 17654.666
$ws.Range("L42").Value = 17654.666
$ws.Range("N42").Value = -18780.666
$ws.Range("H49").Value = 16532.572
$ws.Range("J49").Value = 17654.666
$ws.Range("L49").Value = 17654.666
$ws.Range("N49").Value = -17948.666
$ws.Range("H50").Value = 10941
$ws.Range("I50").Value = 6500
$ws.Range("J50").Value = 12421.333
$ws.Range("K50").Value = 6500
$ws.Range("L50").Value = 12421.333
$ws.Range("M50").Value = -5863
$ws.Range("N50").Value = -13695.333
$ws.Range("H54").Value = 21296
$ws.Range("J54").Value = 21296
$ws.Range("L54").Value = 21296
$ws.Range("N54").Value = -22584
$ws.Range("H132").Value = 4980082
$ws.Range("I132").Value = 1672.8206
$ws.Range("J132").Value = 11914294
$ws.Range("K132").Value = 5018.4618
$ws.Range("L132").Value = 35742882
$ws.Range("M132").Value = -2488.4618
$ws.Range("N132").Value = -35747942

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 625.75
$ws.Range("I113").Value = 596.8889
$ws.Range("K113").Value = 1790.6667
$ws.Range("M113").Value = 379.3332999999998
